$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '27.035.96'
$ws.Range('E2').Value = '  -1.01%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '1.825.97'
$ws.Range('E3').Value = '  -0.39%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '1.008'
$ws.Range('E4').Value = '  -0.41%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '310.22'
$ws.Range('E5').Value = '  -1.57%  '
$ws.Range('E6').Value = '  -0.35%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.4634'
$ws.Range('E7').Value = '  -2.21%  '
$ws.Range('E8').Value = '  -0.91%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.07247'
$ws.Range('E9').Value = '  -2.66%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.8604'
$ws.Range('E10').Value = '  -2.86%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '19.92'
$ws.Range('E11').Value = '  -2.79%  '
$ws.Range('B12').Value = 'TRON'
$ws.Range('C12').Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.07821'
$ws.Range('E12').Value = '  +6.71%  '
$ws.Range('B13').Value = 'WrappedEther'
$ws.Range('C13').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '1.910.33'
$ws.Range('E13').Value = '  +1.56%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '5.332'
$ws.Range('E14').Value = '  -1.91%  '
$ws.Range('B15').Value = 'Chainlink'
$ws.Range('C15').Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '6.509'
$ws.Range('E15').Value = '  -0.85%  '
$ws.Range('B16').Value = 'Litecoin'
$ws.Range('C16').Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '91.82'
$ws.Range('E16').Value = '  -2.31%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '1.008'
$ws.Range('E17').Value = '  -0.28%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '0.000008679'
$ws.Range('E18').Value = '  -1.37%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '1.005'
$ws.Range('E19').Value = '  -0.47%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '14.50'
$ws.Range('E20').Value = '  -1.89%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '26.817.89'
$ws.Range('E21').Value = '  -2.55%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '5.161'
$ws.Range('E22').Value = '  -2.38%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '10.54'
$ws.Range('E23').Value = '  -1.19%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '2.132.56'
$ws.Range('E24').Value = '  +1.76%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '152.08'
$ws.Range('E25').Value = '  +0.00%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '1.839'
$ws.Range('E26').Value = '  -3.04%  '
$ws.Range('E27').Value = '  -2.65%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '2.078'
$ws.Range('E28').Value = '  -3.10%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '5.102'
$ws.Range('E29').Value = '  -2.43%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '115.22'
$ws.Range('E30').Value = '  -1.63%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '0.08815'
$ws.Range('E31').Value = '  -1.96%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '2.956'
$ws.Range('E32').Value = '  +0.37%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '4.433'
$ws.Range('E33').Value = '  -2.45%  '
$ws.Range('E34').Value = '  -3.66%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '0.7212'
$ws.Range('E35').Value = '  -3.87%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '1.080'
$ws.Range('E36').Value = '  -1.22%  '
$ws.Range('E37').Value = '  -2.05%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '2.428'
$ws.Range('E38').Value = '  +1.22%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.01935'
$ws.Range('E39').Value = '  -0.95%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '2.946'
$ws.Range('E40').Value = '  -0.88%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '7.183'
$ws.Range('E41').Value = '  -0.64%  '
$ws.Range('E42').Value = '  -2.60%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.1628'
$ws.Range('E43').Value = '  -1.94%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.8583'
$ws.Range('E44').Value = '  -15.14%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '8.180'
$ws.Range('E45').Value = '  -3.51%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.4801'
$ws.Range('E46').Value = '  -2.76%  '
$ws.Range('E47').Value = '  -0.44%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '10.13'
$ws.Range('E48').Value = '  -3.72%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '102.63'
$ws.Range('E49').Value = '  -2.28%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.06243'
$ws.Range('E50').Value = '  -0.90%  '
$ws.Range('E51').Value = '  -3.27%  '
